$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells are treated as text so values like "242.30" or
# "0.000008525" are not auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.324.89'
$ws.Range("E2").Value = '  -0.13%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.877.78'
$ws.Range("E3").Value = '  +0.18%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7107'
$ws.Range("E5").Value = '  -0.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.30'
$ws.Range("E6").Value = '  +0.17%  '

$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.08005'
$ws.Range("E8").Value = '  +3.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3164'
$ws.Range("E9").Value = '  +1.70%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.98'
$ws.Range("E10").Value = '  -0.56%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08304'
$ws.Range("E11").Value = '  -1.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.876.92'
$ws.Range("E12").Value = '  -0.39%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.252'
$ws.Range("E13").Value = '  -0.13%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.52'
$ws.Range("E14").Value = '  +3.61%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7150'
$ws.Range("E15").Value = '  +0.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.377'
$ws.Range("E16").Value = '  +4.72%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008525'
$ws.Range("E17").Value = '  +3.53%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.346.76'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.27'
$ws.Range("E19").Value = '  +1.08%  '

$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.138.86'
$ws.Range("E20").Value = '  +0.69%  '

$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.29'
$ws.Range("E21").Value = '  +0.38%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.806'
$ws.Range("E23").Value = '  +0.19%  '

$ws.Range("E24").Value = '  +0.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1559'
$ws.Range("E25").Value = '  -2.24%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.072'
$ws.Range("E26").Value = '  +0.22%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.71'
$ws.Range("E27").Value = '  -0.28%  '

$ws.Range("E28").Value = '  +0.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.508'
$ws.Range("E29").Value = '  -0.31%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.423'
$ws.Range("E30").Value = '  +0.00%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.324'
$ws.Range("E31").Value = '  -0.14%  '

$ws.Range("E32").Value = '  -7.15%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05394'
$ws.Range("E33").Value = '  +1.54%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.938'
$ws.Range("E34").Value = '  +0.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7711'
$ws.Range("E35").Value = '  +4.20%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.183'
$ws.Range("E36").Value = '  +0.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.687'
$ws.Range("E37").Value = '  -0.43%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01886'
$ws.Range("E38").Value = '  +0.76%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.262.23'
$ws.Range("E39").Value = '  +2.53%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.752'
$ws.Range("E40").Value = '  +0.75%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.484'
$ws.Range("E41").Value = '  -0.68%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '113.18'
$ws.Range("E42").Value = '  +2.42%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9061'
$ws.Range("E43").Value = '  +1.83%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '74.17'
$ws.Range("E44").Value = '  +1.54%  '

$ws.Range("E45").Value = '  +7.70%  '

$ws.Range("E46").Value = '  +0.10%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.030.24'
$ws.Range("E47").Value = '  +0.43%  '

$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5227'
$ws.Range("E48").Value = '  +0.23%  '

$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.802'
$ws.Range("E49").Value = '  -0.43%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.462'
$ws.Range("E50").Value = '  +0.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4364'
$ws.Range("E51").Value = '  +1.06%  '
